$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 9960.091
$ws.Cells.Item(9, 9).Value = 12757.625
$ws.Cells.Item(9, 10).Value = 2500
$ws.Cells.Item(9, 11).Value = 12757.625
$ws.Cells.Item(9, 12).Value = 2500
$ws.Cells.Item(9, 13).Value = -12588.625
$ws.Cells.Item(9, 14).Value = -2838

$ws.Cells.Item(69, 8).Value = 7005
$ws.Cells.Item(69, 9).Value = 7005
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 21015
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = -20141
$ws.Cells.Item(69, 14).ClearContents()

$ws.Cells.Item(70, 8).Value = 3613
$ws.Cells.Item(70, 9).Value = 3000
$ws.Cells.Item(70, 11).Value = 9000
$ws.Cells.Item(70, 13).Value = -8730

$ws.Cells.Item(72, 8).Value = 7005
$ws.Cells.Item(72, 9).Value = 7005
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 11).Value = 63045
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 13).Value = -58677
$ws.Cells.Item(72, 14).ClearContents()

$ws.Cells.Item(73, 8).Value = 3613
$ws.Cells.Item(73, 9).Value = 3000
$ws.Cells.Item(73, 11).Value = 9000
$ws.Cells.Item(73, 13).Value = -8064

$ws.Cells.Item(74, 8).Value = 4928.4116
$ws.Cells.Item(74, 9).Value = 3444.8333
$ws.Cells.Item(74, 11).Value = 3444.8333
$ws.Cells.Item(74, 13).Value = -2508.8333

$ws.Cells.Item(76, 8).Value = 9769.611000000001
$ws.Cells.Item(76, 9).Value = 10579.077
$ws.Cells.Item(76, 10).Value = 7665
$ws.Cells.Item(76, 11).Value = 10579.077
$ws.Cells.Item(76, 12).Value = 7665
$ws.Cells.Item(76, 13).Value = -10264.077
$ws.Cells.Item(76, 14).Value = -8295

$ws.Cells.Item(77, 8).Value = 4928.4116
$ws.Cells.Item(77, 9).Value = 3444.8333
$ws.Cells.Item(77, 11).Value = 17224.1665
$ws.Cells.Item(77, 13).Value = -12544.1665

$ws.Cells.Item(79, 8).Value = 9769.611000000001
$ws.Cells.Item(79, 9).Value = 10579.077
$ws.Cells.Item(79, 10).Value = 7665
$ws.Cells.Item(79, 11).Value = 10579.077
$ws.Cells.Item(79, 12).Value = 7665
$ws.Cells.Item(79, 13).Value = -9487.076999999999
$ws.Cells.Item(79, 14).Value = -9849

$ws.Cells.Item(80, 8).Value = 605.2143
$ws.Cells.Item(80, 9).Value = 387.33334
$ws.Cells.Item(80, 10).Value = 856.61536
$ws.Cells.Item(80, 11).Value = 1162.00002
$ws.Cells.Item(80, 12).Value = 2569.84608
$ws.Cells.Item(80, 13).Value = -164.0000199999999
$ws.Cells.Item(80, 14).Value = -4565.84608

$ws.Cells.Item(83, 8).Value = 605.2143
$ws.Cells.Item(83, 9).Value = 387.33334
$ws.Cells.Item(83, 10).Value = 856.61536
$ws.Cells.Item(83, 11).Value = 3486.00006
$ws.Cells.Item(83, 12).Value = 7709.53824
$ws.Cells.Item(83, 13).Value = 1505.99994
$ws.Cells.Item(83, 14).Value = -17693.53824

$ws.Cells.Item(138, 8).Value = 7936.7905
$ws.Cells.Item(138, 9).Value = 3603
$ws.Cells.Item(138, 10).Value = 8927.370999999999
$ws.Cells.Item(138, 11).Value = 10809
$ws.Cells.Item(138, 12).Value = 26782.113
$ws.Cells.Item(138, 13).Value = -5669
$ws.Cells.Item(138, 14).Value = -37062.113

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 17866390
$ws.Cells.Item(32, 9).Value = 21282634
$ws.Cells.Item(32, 11).Value = 21282634
$ws.Cells.Item(32, 13).Value = -21282347

$ws.Cells.Item(61, 8).Value = 33336570
$ws.Cells.Item(61, 9).Value = 40001784
$ws.Cells.Item(61, 11).Value = 40001784
$ws.Cells.Item(61, 13).Value = -40001572

$ws.Cells.Item(102, 8).Value = 2829.2666
$ws.Cells.Item(102, 9).Value = 2674.2144
$ws.Cells.Item(102, 10).Value = 5000
$ws.Cells.Item(102, 11).Value = 2674.2144
$ws.Cells.Item(102, 12).Value = 5000
$ws.Cells.Item(102, 13).Value = -1052.2144
$ws.Cells.Item(102, 14).Value = -8244

$ws.Cells.Item(136, 8).Value = 33336570
$ws.Cells.Item(136, 9).Value = 40001784
$ws.Cells.Item(136, 11).Value = 120005352
$ws.Cells.Item(136, 13).Value = -120002802

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(25, 8).Value = 2975.5715
$ws.Cells.Item(25, 9).Value = 3309.3333
$ws.Cells.Item(25, 11).Value = 3309.3333
$ws.Cells.Item(25, 13).Value = -3074.3333

$ws.Cells.Item(37, 8).Value = 2073.2273
$ws.Cells.Item(37, 9).Value = 1408.2727
$ws.Cells.Item(37, 10).Value = 2738.182
$ws.Cells.Item(37, 11).Value = 1408.2727
$ws.Cells.Item(37, 12).Value = 2738.182
$ws.Cells.Item(37, 13).Value = -1271.2727
$ws.Cells.Item(37, 14).Value = -3012.182

$ws.Cells.Item(80, 8).Value = 8142.7144
$ws.Cells.Item(80, 9).Value = 8250
$ws.Cells.Item(80, 10).Value = 7999.6665
$ws.Cells.Item(80, 11).Value = 8250
$ws.Cells.Item(80, 12).Value = 7999.6665
$ws.Cells.Item(80, 13).Value = -7252
$ws.Cells.Item(80, 14).Value = -9995.666499999999

$ws.Cells.Item(83, 8).Value = 8142.7144
$ws.Cells.Item(83, 9).Value = 8250
$ws.Cells.Item(83, 10).Value = 7999.6665
$ws.Cells.Item(83, 11).Value = 41250
$ws.Cells.Item(83, 12).Value = 39998.3325
$ws.Cells.Item(83, 13).Value = -36258
$ws.Cells.Item(83, 14).Value = -49982.3325

$ws.Cells.Item(86, 8).Value = 22510.883
$ws.Cells.Item(86, 9).Value = 12928.5
$ws.Cells.Item(86, 11).Value = 12928.5
$ws.Cells.Item(86, 13).Value = -11805.5

$ws.Cells.Item(89, 8).Value = 22510.883
$ws.Cells.Item(89, 9).Value = 12928.5
$ws.Cells.Item(89, 11).Value = 64642.5
$ws.Cells.Item(89, 13).Value = -59026.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1728.1111
$ws.Cells.Item(58, 9).Value = 1662.875
$ws.Cells.Item(58, 11).Value = 1662.875
$ws.Cells.Item(58, 13).Value = -1459.875

$ws.Cells.Item(86, 8).Value = 5772.875
$ws.Cells.Item(86, 9).Value = 4703
$ws.Cells.Item(86, 10).Value = 6129.5
$ws.Cells.Item(86, 11).Value = 4703
$ws.Cells.Item(86, 12).Value = 6129.5
$ws.Cells.Item(86, 13).Value = -3580
$ws.Cells.Item(86, 14).Value = -8375.5

$ws.Cells.Item(89, 8).Value = 5772.875
$ws.Cells.Item(89, 9).Value = 4703
$ws.Cells.Item(89, 10).Value = 6129.5
$ws.Cells.Item(89, 11).Value = 23515
$ws.Cells.Item(89, 12).Value = 30647.5
$ws.Cells.Item(89, 13).Value = -17899
$ws.Cells.Item(89, 14).Value = -41879.5

$ws.Cells.Item(132, 8).Value = 3816.0625
$ws.Cells.Item(132, 9).Value = 2907.75
$ws.Cells.Item(132, 11).Value = 8723.25
$ws.Cells.Item(132, 13).Value = -6193.25

$ws.Cells.Item(136, 8).Value = 1728.1111
$ws.Cells.Item(136, 9).Value = 1662.875
$ws.Cells.Item(136, 11).Value = 4988.625
$ws.Cells.Item(136, 13).Value = -2438.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 2381021.2
$ws.Cells.Item(2, 9).Value = 3571494.5
$ws.Cells.Item(2, 11).Value = 3571494.5
$ws.Cells.Item(2, 13).Value = -3571381.5

$ws.Cells.Item(12, 8).Value = 10005000
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 13).ClearContents()

$ws.Cells.Item(14, 8).Value = 10169135
$ws.Cells.Item(14, 9).Value = 6461419
$ws.Cells.Item(14, 11).Value = 6461419
$ws.Cells.Item(14, 13).Value = -6461251

$ws.Cells.Item(70, 8).Value = 5058.5864
$ws.Cells.Item(70, 9).Value = 4768.273
$ws.Cells.Item(70, 11).Value = 4768.273
$ws.Cells.Item(70, 13).Value = -4498.273

$ws.Cells.Item(73, 8).Value = 5058.5864
$ws.Cells.Item(73, 9).Value = 4768.273
$ws.Cells.Item(73, 11).Value = 4768.273
$ws.Cells.Item(73, 13).Value = -3832.273

$ws.Cells.Item(80, 8).Value = 4891
$ws.Cells.Item(80, 9).Value = 4251.25
$ws.Cells.Item(80, 10).Value = 5210.875
$ws.Cells.Item(80, 11).Value = 4251.25
$ws.Cells.Item(80, 12).Value = 5210.875
$ws.Cells.Item(80, 13).Value = -3253.25
$ws.Cells.Item(80, 14).Value = -7206.875

$ws.Cells.Item(83, 8).Value = 4891
$ws.Cells.Item(83, 9).Value = 4251.25
$ws.Cells.Item(83, 10).Value = 5210.875
$ws.Cells.Item(83, 11).Value = 21256.25
$ws.Cells.Item(83, 12).Value = 26054.375
$ws.Cells.Item(83, 13).Value = -16264.25
$ws.Cells.Item(83, 14).Value = -36038.375

$ws.Cells.Item(113, 8).Value = 7471.1665
$ws.Cells.Item(113, 9).Value = 8999
$ws.Cells.Item(113, 10).Value = 6707.25
$ws.Cells.Item(113, 11).Value = 8999
$ws.Cells.Item(113, 12).Value = 6707.25
$ws.Cells.Item(113, 13).Value = -6829
$ws.Cells.Item(113, 14).Value = -11047.25

$ws.Cells.Item(135, 8).Value = 86729.64999999999
$ws.Cells.Item(135, 10).Value = 86729.64999999999
$ws.Cells.Item(135, 12).Value = 86729.64999999999
$ws.Cells.Item(135, 14).Value = -96869.64999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4445.8086
$ws.Cells.Item(7, 9).Value = 4149
$ws.Cells.Item(7, 10).Value = 4730.25
$ws.Cells.Item(7, 11).Value = 4149
$ws.Cells.Item(7, 12).Value = 4730.25
$ws.Cells.Item(7, 13).Value = -4037
$ws.Cells.Item(7, 14).Value = -4954.25

$ws.Cells.Item(55, 8).Value = 637.7222
$ws.Cells.Item(55, 9).Value = 332.9
$ws.Cells.Item(55, 11).Value = 332.9
$ws.Cells.Item(55, 13).Value = -159.9

$ws.Cells.Item(104, 8).Value = 10955.5
$ws.Cells.Item(104, 10).Value = 10955.5
$ws.Cells.Item(104, 12).Value = 10955.5
$ws.Cells.Item(104, 14).Value = -17943.5

$ws.Cells.Item(126, 8).Value = 4445.8086
$ws.Cells.Item(126, 9).Value = 4149
$ws.Cells.Item(126, 10).Value = 4730.25
$ws.Cells.Item(126, 11).Value = 12447
$ws.Cells.Item(126, 12).Value = 14190.75
$ws.Cells.Item(126, 13).Value = -9977
$ws.Cells.Item(126, 14).Value = -19130.75

$ws.Cells.Item(136, 8).Value = 4378.8623
$ws.Cells.Item(136, 9).Value = 4135.2085
$ws.Cells.Item(136, 11).Value = 12405.6255
$ws.Cells.Item(136, 13).Value = -9855.625499999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).ClearContents()

$ws.Cells.Item(132, 8).Value = 5854.4287
$ws.Cells.Item(132, 9).Value = 5967.205
$ws.Cells.Item(132, 11).Value = 17901.615
$ws.Cells.Item(132, 13).Value = -15371.615

$ws.Cells.Item(136, 8).Value = 1746.3
$ws.Cells.Item(136, 9).Value = 1084.5
$ws.Cells.Item(136, 11).Value = 3253.5
$ws.Cells.Item(136, 13).Value = -703.5
